$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24
$ws.Cells.Item(24, 2).Value = 6281792
$ws.Cells.Item(24, 6).Value = "Ituano"
$ws.Cells.Item(24, 7).Value = "Atletico GO"
$ws.Cells.Item(24, 9).Value = 1
$ws.Cells.Item(24, 10).Value = "D"
$ws.Cells.Item(24, 11).Value = 2.45
$ws.Cells.Item(24, 12).Value = 3.2
$ws.Cells.Item(24, 13).Value = 2.9
$ws.Cells.Item(24, 14).Value = 2.25
$ws.Cells.Item(24, 15).Value = 3.25
$ws.Cells.Item(24, 16).Value = 3.5
$ws.Cells.Item(24, 17).Value = -0.25
$ws.Cells.Item(24, 18).Value = 1.925
$ws.Cells.Item(24, 19).Value = 1.925
$ws.Cells.Item(24, 21).Value = 1.925
$ws.Cells.Item(24, 22).Value = 1.925
$ws.Cells.Item(24, 23).Value = -1
$ws.Cells.Item(24, 24).Value = 2.25
$ws.Cells.Item(24, 26).Value = -0.5
$ws.Cells.Item(24, 27).Value = 0.4625
$ws.Cells.Item(24, 28).Value = -0.5
$ws.Cells.Item(24, 29).Value = 0.4625

# Row 25
$ws.Cells.Item(25, 2).Value = 6281967
$ws.Cells.Item(25, 6).Value = "Gremio Novorizontino"
$ws.Cells.Item(25, 7).Value = "Sampaio Correa"
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = "H"
$ws.Cells.Item(25, 11).Value = 1.666
$ws.Cells.Item(25, 12).Value = 3.4
$ws.Cells.Item(25, 13).Value = 6
$ws.Cells.Item(25, 14).Value = 1.615
$ws.Cells.Item(25, 15).Value = 3.6
$ws.Cells.Item(25, 16).Value = 6
$ws.Cells.Item(25, 17).Value = -0.75
$ws.Cells.Item(25, 18).Value = 1.775
$ws.Cells.Item(25, 19).Value = 2.025
$ws.Cells.Item(25, 21).Value = 1.975
$ws.Cells.Item(25, 22).Value = 1.825
$ws.Cells.Item(25, 23).Value = 0.615
$ws.Cells.Item(25, 24).Value = -1
$ws.Cells.Item(25, 26).Value = 0.3875
$ws.Cells.Item(25, 27).Value = -0.5
$ws.Cells.Item(25, 28).Value = -1
$ws.Cells.Item(25, 29).Value = 0.825

# Row 27
$ws.Cells.Item(27, 2).Value = 6285527
$ws.Cells.Item(27, 6).Value = "Ponte Preta"
$ws.Cells.Item(27, 7).Value = "Sport Recife"
$ws.Cells.Item(27, 9).Value = 1
$ws.Cells.Item(27, 10).Value = "D"
$ws.Cells.Item(27, 11).Value = 3
$ws.Cells.Item(27, 12).Value = 3
$ws.Cells.Item(27, 13).Value = 2.55
$ws.Cells.Item(27, 14).Value = 3.6
$ws.Cells.Item(27, 15).Value = 3
$ws.Cells.Item(27, 19).Value = 1.875
$ws.Cells.Item(27, 21).Value = 1.8
$ws.Cells.Item(27, 24).Value = 2
$ws.Cells.Item(27, 25).Value = -1
$ws.Cells.Item(27, 26).Value = 0.4875
$ws.Cells.Item(27, 27).Value = -0.5
$ws.Cells.Item(27, 28).Value = 0
$ws.Cells.Item(27, 29).Value = -0

# Row 28
$ws.Cells.Item(28, 2).Value = 6289131
$ws.Cells.Item(28, 6).Value = "Londrina"
$ws.Cells.Item(28, 7).Value = "Mirassol"
$ws.Cells.Item(28, 9).Value = 2
$ws.Cells.Item(28, 10).Value = "A"
$ws.Cells.Item(28, 11).Value = 2.75
$ws.Cells.Item(28, 12).Value = 3.1
$ws.Cells.Item(28, 13).Value = 2.7
$ws.Cells.Item(28, 14).Value = 3.4
$ws.Cells.Item(28, 15).Value = 3.2
$ws.Cells.Item(28, 19).Value = 1.825
$ws.Cells.Item(28, 21).Value = 1.75
$ws.Cells.Item(28, 24).Value = -1
$ws.Cells.Item(28, 25).Value = 1.2
$ws.Cells.Item(28, 26).Value = -1
$ws.Cells.Item(28, 27).Value = 0.825
$ws.Cells.Item(28, 28).Value = 0.75
$ws.Cells.Item(28, 29).Value = -1

# Row 47
$ws.Cells.Item(47, 2).Value = 6285666
$ws.Cells.Item(47, 6).Value = "ABC"
$ws.Cells.Item(47, 7).Value = "Atletico GO"
$ws.Cells.Item(47, 8).Value = 1
$ws.Cells.Item(47, 9).Value = 1
$ws.Cells.Item(47, 11).Value = 2.75
$ws.Cells.Item(47, 12).Value = 3.1
$ws.Cells.Item(47, 13).Value = 2.45
$ws.Cells.Item(47, 14).Value = 3
$ws.Cells.Item(47, 15).Value = 3
$ws.Cells.Item(47, 16).Value = 2.5
$ws.Cells.Item(47, 17).Value = 0
$ws.Cells.Item(47, 18).Value = 2.1
$ws.Cells.Item(47, 19).Value = 1.775
$ws.Cells.Item(47, 20).Value = 2
$ws.Cells.Item(47, 21).Value = 1.975
$ws.Cells.Item(47, 22).Value = 1.875
$ws.Cells.Item(47, 24).Value = 2
$ws.Cells.Item(47, 26).Value = 0
$ws.Cells.Item(47, 27).Value = -0
$ws.Cells.Item(47, 28).Value = 0
$ws.Cells.Item(47, 29).Value = -0

# Row 48
$ws.Cells.Item(48, 2).Value = 6282100
$ws.Cells.Item(48, 6).Value = "Ceara"
$ws.Cells.Item(48, 7).Value = "Avai"
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 11).Value = 1.666
$ws.Cells.Item(48, 12).Value = 3.5
$ws.Cells.Item(48, 13).Value = 4.75
$ws.Cells.Item(48, 14).Value = 1.75
$ws.Cells.Item(48, 15).Value = 3.6
$ws.Cells.Item(48, 16).Value = 4.75
$ws.Cells.Item(48, 17).Value = -0.75
$ws.Cells.Item(48, 18).Value = 2
$ws.Cells.Item(48, 19).Value = 1.8
$ws.Cells.Item(48, 20).Value = 2.25
$ws.Cells.Item(48, 21).Value = 1.9
$ws.Cells.Item(48, 22).Value = 1.9
$ws.Cells.Item(48, 24).Value = 2.6
$ws.Cells.Item(48, 26).Value = -1
$ws.Cells.Item(48, 27).Value = 0.8
$ws.Cells.Item(48, 28).Value = -1
$ws.Cells.Item(48, 29).Value = 0.8999999999999999

# Row 100
$ws.Cells.Item(100, 2).Value = 6286300
$ws.Cells.Item(100, 6).Value = "Vitoria"
$ws.Cells.Item(100, 7).Value = "Chapecoense"
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = "H"
$ws.Cells.Item(100, 11).Value = 1.7
$ws.Cells.Item(100, 12).Value = 3.25
$ws.Cells.Item(100, 13).Value = 4.5
$ws.Cells.Item(100, 14).Value = 1.7
$ws.Cells.Item(100, 15).Value = 3.5
$ws.Cells.Item(100, 16).Value = 5.75
$ws.Cells.Item(100, 17).Value = -0.75
$ws.Cells.Item(100, 18).Value = 1.9
$ws.Cells.Item(100, 19).Value = 1.95
$ws.Cells.Item(100, 20).Value = 2.25
$ws.Cells.Item(100, 21).Value = 2.05
$ws.Cells.Item(100, 22).Value = 1.8
$ws.Cells.Item(100, 23).Value = 0.7
$ws.Cells.Item(100, 25).Value = -1
$ws.Cells.Item(100, 26).Value = 0.45
$ws.Cells.Item(100, 27).Value = -0.5
$ws.Cells.Item(100, 28).Value = -1
$ws.Cells.Item(100, 29).Value = 0.8

# Row 101
$ws.Cells.Item(101, 2).Value = 6281811
$ws.Cells.Item(101, 6).Value = "Londrina"
$ws.Cells.Item(101, 7).Value = "Botafogo SP"
$ws.Cells.Item(101, 9).Value = 2
$ws.Cells.Item(101, 10).Value = "A"
$ws.Cells.Item(101, 11).Value = 2.5
$ws.Cells.Item(101, 12).Value = 3
$ws.Cells.Item(101, 13).Value = 2.625
$ws.Cells.Item(101, 14).Value = 2.6
$ws.Cells.Item(101, 15).Value = 3
$ws.Cells.Item(101, 16).Value = 2.9
$ws.Cells.Item(101, 17).Value = 0
$ws.Cells.Item(101, 18).Value = 1.775
$ws.Cells.Item(101, 19).Value = 2.025
$ws.Cells.Item(101, 20).Value = 1.75
$ws.Cells.Item(101, 21).Value = 1.8
$ws.Cells.Item(101, 22).Value = 2
$ws.Cells.Item(101, 23).Value = -1
$ws.Cells.Item(101, 25).Value = 1.9
$ws.Cells.Item(101, 26).Value = -1
$ws.Cells.Item(101, 27).Value = 1.025
$ws.Cells.Item(101, 28).Value = 0.8
$ws.Cells.Item(101, 29).Value = -1

# Row 112
$ws.Cells.Item(112, 2).Value = 6281874
$ws.Cells.Item(112, 6).Value = "Sampaio Correa"
$ws.Cells.Item(112, 7).Value = "Botafogo SP"
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 11).Value = 2.3
$ws.Cells.Item(112, 12).Value = 2.8
$ws.Cells.Item(112, 13).Value = 3.1
$ws.Cells.Item(112, 14).Value = 2
$ws.Cells.Item(112, 15).Value = 3.1
$ws.Cells.Item(112, 16).Value = 4.333
$ws.Cells.Item(112, 17).Value = -0.5
$ws.Cells.Item(112, 18).Value = 2.025
$ws.Cells.Item(112, 19).Value = 1.825
$ws.Cells.Item(112, 20).Value = 1.75
$ws.Cells.Item(112, 21).Value = 1.825
$ws.Cells.Item(112, 22).Value = 2.025
$ws.Cells.Item(112, 24).Value = 2.1
$ws.Cells.Item(112, 27).Value = 0.825
$ws.Cells.Item(112, 28).Value = -1
$ws.Cells.Item(112, 29).Value = 1.025

# Row 113
$ws.Cells.Item(113, 2).Value = 6289120
$ws.Cells.Item(113, 6).Value = "Mirassol"
$ws.Cells.Item(113, 7).Value = "Avai"
$ws.Cells.Item(113, 8).Value = 2
$ws.Cells.Item(113, 9).Value = 2
$ws.Cells.Item(113, 11).Value = 1.615
$ws.Cells.Item(113, 12).Value = 3.4
$ws.Cells.Item(113, 13).Value = 5
$ws.Cells.Item(113, 14).Value = 1.615
$ws.Cells.Item(113, 15).Value = 3.6
$ws.Cells.Item(113, 16).Value = 6
$ws.Cells.Item(113, 17).Value = -0.75
$ws.Cells.Item(113, 18).Value = 1.825
$ws.Cells.Item(113, 19).Value = 1.975
$ws.Cells.Item(113, 20).Value = 2
$ws.Cells.Item(113, 21).Value = 1.8
$ws.Cells.Item(113, 22).Value = 2
$ws.Cells.Item(113, 24).Value = 2.6
$ws.Cells.Item(113, 27).Value = 0.9750000000000001
$ws.Cells.Item(113, 28).Value = 0.8
$ws.Cells.Item(113, 29).Value = -1

# Row 115
$ws.Cells.Item(115, 2).Value = 6285545
$ws.Cells.Item(115, 6).Value = "Criciuma"
$ws.Cells.Item(115, 7).Value = "Ponte Preta"
$ws.Cells.Item(115, 8).Value = 2
$ws.Cells.Item(115, 9).Value = 1
$ws.Cells.Item(115, 11).Value = 1.727
$ws.Cells.Item(115, 12).Value = 3.2
$ws.Cells.Item(115, 13).Value = 4.5
$ws.Cells.Item(115, 14).Value = 1.75
$ws.Cells.Item(115, 15).Value = 3.5
$ws.Cells.Item(115, 16).Value = 5.25
$ws.Cells.Item(115, 17).Value = -0.75
$ws.Cells.Item(115, 18).Value = 2
$ws.Cells.Item(115, 19).Value = 1.8
$ws.Cells.Item(115, 21).Value = 1.875
$ws.Cells.Item(115, 22).Value = 1.925
$ws.Cells.Item(115, 23).Value = 0.75
$ws.Cells.Item(115, 26).Value = 0.5
$ws.Cells.Item(115, 27).Value = -0.5
$ws.Cells.Item(115, 28).Value = 0.875
$ws.Cells.Item(115, 29).Value = -1

# Row 116
$ws.Cells.Item(116, 2).Value = 6281816
$ws.Cells.Item(116, 6).Value = "Ituano"
$ws.Cells.Item(116, 7).Value = "Tombense MG"
$ws.Cells.Item(116, 8).Value = 1
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 11).Value = 2.1
$ws.Cells.Item(116, 12).Value = 3
$ws.Cells.Item(116, 13).Value = 3.25
$ws.Cells.Item(116, 14).Value = 1.95
$ws.Cells.Item(116, 15).Value = 3.2
$ws.Cells.Item(116, 16).Value = 4.5
$ws.Cells.Item(116, 17).Value = -0.5
$ws.Cells.Item(116, 18).Value = 1.975
$ws.Cells.Item(116, 19).Value = 1.825
$ws.Cells.Item(116, 21).Value = 1.95
$ws.Cells.Item(116, 22).Value = 1.85
$ws.Cells.Item(116, 23).Value = 0.95
$ws.Cells.Item(116, 26).Value = 0.9750000000000001
$ws.Cells.Item(116, 27).Value = -1
$ws.Cells.Item(116, 28).Value = -1
$ws.Cells.Item(116, 29).Value = 0.8500000000000001

# Row 117
$ws.Cells.Item(117, 2).Value = 6285671
$ws.Cells.Item(117, 6).Value = "Vitoria"
$ws.Cells.Item(117, 7).Value = "ABC"
$ws.Cells.Item(117, 8).Value = 2
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = "H"
$ws.Cells.Item(117, 11).Value = 1.4
$ws.Cells.Item(117, 12).Value = 4
$ws.Cells.Item(117, 13).Value = 6.5
$ws.Cells.Item(117, 14).Value = 1.5
$ws.Cells.Item(117, 15).Value = 4
$ws.Cells.Item(117, 16).Value = 7.5
$ws.Cells.Item(117, 17).Value = -1
$ws.Cells.Item(117, 18).Value = 1.875
$ws.Cells.Item(117, 19).Value = 1.925
$ws.Cells.Item(117, 20).Value = 2
$ws.Cells.Item(117, 21).Value = 1.85
$ws.Cells.Item(117, 22).Value = 1.95
$ws.Cells.Item(117, 23).Value = 0.5
$ws.Cells.Item(117, 25).Value = -1
$ws.Cells.Item(117, 26).Value = 0.875
$ws.Cells.Item(117, 27).Value = -1
$ws.Cells.Item(117, 28).Value = 0
$ws.Cells.Item(117, 29).Value = -0

# Row 118
$ws.Cells.Item(118, 2).Value = 6281815
$ws.Cells.Item(118, 6).Value = "Guarani"
$ws.Cells.Item(118, 7).Value = "Ceara"
$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 11).Value = 1.95
$ws.Cells.Item(118, 12).Value = 3.25
$ws.Cells.Item(118, 13).Value = 3.4
$ws.Cells.Item(118, 14).Value = 2.05
$ws.Cells.Item(118, 15).Value = 3.25
$ws.Cells.Item(118, 16).Value = 4
$ws.Cells.Item(118, 17).Value = -0.5
$ws.Cells.Item(118, 18).Value = 2.025
$ws.Cells.Item(118, 21).Value = 1.95
$ws.Cells.Item(118, 22).Value = 1.85
$ws.Cells.Item(118, 24).Value = 2.25
$ws.Cells.Item(118, 26).Value = -1
$ws.Cells.Item(118, 27).Value = 0.7749999999999999
$ws.Cells.Item(118, 28).Value = -1
$ws.Cells.Item(118, 29).Value = 0.8500000000000001

# Row 120
$ws.Cells.Item(120, 2).Value = 6287040
$ws.Cells.Item(120, 6).Value = "Londrina"
$ws.Cells.Item(120, 7).Value = "Chapecoense"
$ws.Cells.Item(120, 8).Value = 1
$ws.Cells.Item(120, 9).Value = 1
$ws.Cells.Item(120, 11).Value = 2.4
$ws.Cells.Item(120, 12).Value = 2.875
$ws.Cells.Item(120, 13).Value = 2.875
$ws.Cells.Item(120, 14).Value = 2.5
$ws.Cells.Item(120, 15).Value = 2.875
$ws.Cells.Item(120, 16).Value = 3.2
$ws.Cells.Item(120, 17).Value = -0.25
$ws.Cells.Item(120, 18).Value = 2.1
$ws.Cells.Item(120, 21).Value = 2.05
$ws.Cells.Item(120, 22).Value = 1.8
$ws.Cells.Item(120, 24).Value = 1.875
$ws.Cells.Item(120, 26).Value = -0.5
$ws.Cells.Item(120, 27).Value = 0.3875
$ws.Cells.Item(120, 28).Value = 0
$ws.Cells.Item(120, 29).Value = -0

# Row 121
$ws.Cells.Item(121, 2).Value = 6285544
$ws.Cells.Item(121, 6).Value = "Vila Nova"
$ws.Cells.Item(121, 7).Value = "Sport Recife"
$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 9).Value = 1
$ws.Cells.Item(121, 10).Value = "A"
$ws.Cells.Item(121, 11).Value = 2.2
$ws.Cells.Item(121, 12).Value = 2.875
$ws.Cells.Item(121, 13).Value = 3.2
$ws.Cells.Item(121, 14).Value = 2.25
$ws.Cells.Item(121, 15).Value = 2.875
$ws.Cells.Item(121, 16).Value = 3.8
$ws.Cells.Item(121, 17).Value = -0.25
$ws.Cells.Item(121, 18).Value = 1.9
$ws.Cells.Item(121, 19).Value = 1.95
$ws.Cells.Item(121, 20).Value = 1.75
$ws.Cells.Item(121, 21).Value = 1.825
$ws.Cells.Item(121, 22).Value = 2.025
$ws.Cells.Item(121, 23).Value = -1
$ws.Cells.Item(121, 25).Value = 2.8
$ws.Cells.Item(121, 26).Value = -1
$ws.Cells.Item(121, 27).Value = 0.95
$ws.Cells.Item(121, 28).Value = -1
$ws.Cells.Item(121, 29).Value = 1.025

# Row 130
$ws.Cells.Item(130, 2).Value = 6959080
$ws.Cells.Item(130, 6).Value = "Atletico GO"
$ws.Cells.Item(130, 7).Value = "Tombense MG"
$ws.Cells.Item(130, 8).Value = 3
$ws.Cells.Item(130, 9).Value = 2
$ws.Cells.Item(130, 11).Value = 1.7
$ws.Cells.Item(130, 12).Value = 3.25
$ws.Cells.Item(130, 13).Value = 4.5
$ws.Cells.Item(130, 14).Value = 1.727
$ws.Cells.Item(130, 15).Value = 3.6
$ws.Cells.Item(130, 16).Value = 5.25
$ws.Cells.Item(130, 17).Value = -0.75
$ws.Cells.Item(130, 18).Value = 1.975
$ws.Cells.Item(130, 19).Value = 1.825
$ws.Cells.Item(130, 20).Value = 2.25
$ws.Cells.Item(130, 21).Value = 1.975
$ws.Cells.Item(130, 22).Value = 1.825
$ws.Cells.Item(130, 23).Value = 0.7270000000000001
$ws.Cells.Item(130, 26).Value = 0.4875
$ws.Cells.Item(130, 27).Value = -0.5
$ws.Cells.Item(130, 28).Value = 0.9750000000000001
$ws.Cells.Item(130, 29).Value = -1

# Row 131
$ws.Cells.Item(131, 2).Value = 6286301
$ws.Cells.Item(131, 6).Value = "Londrina"
$ws.Cells.Item(131, 7).Value = "Vitoria"
$ws.Cells.Item(131, 8).Value = 2
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 11).Value = 3.1
$ws.Cells.Item(131, 12).Value = 3.1
$ws.Cells.Item(131, 13).Value = 2.15
$ws.Cells.Item(131, 14).Value = 3
$ws.Cells.Item(131, 15).Value = 3.2
$ws.Cells.Item(131, 16).Value = 2.45
$ws.Cells.Item(131, 17).Value = 0.25
$ws.Cells.Item(131, 18).Value = 1.725
$ws.Cells.Item(131, 19).Value = 2.075
$ws.Cells.Item(131, 20).Value = 1.75
$ws.Cells.Item(131, 21).Value = 1.775
$ws.Cells.Item(131, 22).Value = 2.025
$ws.Cells.Item(131, 23).Value = 2
$ws.Cells.Item(131, 26).Value = 0.7250000000000001
$ws.Cells.Item(131, 27).Value = -1
$ws.Cells.Item(131, 28).Value = 0.3875
$ws.Cells.Item(131, 29).Value = -0.5

# Row 196
$ws.Cells.Item(196, 2).Value = 6285766
$ws.Cells.Item(196, 6).Value = "Atletico GO"
$ws.Cells.Item(196, 7).Value = "Criciuma"
$ws.Cells.Item(196, 8).Value = 3
$ws.Cells.Item(196, 9).Value = 1
$ws.Cells.Item(196, 10).Value = "H"
$ws.Cells.Item(196, 11).Value = 1.909
$ws.Cells.Item(196, 12).Value = 3.1
$ws.Cells.Item(196, 13).Value = 4
$ws.Cells.Item(196, 14).Value = 1.909
$ws.Cells.Item(196, 15).Value = 3.25
$ws.Cells.Item(196, 16).Value = 4.5
$ws.Cells.Item(196, 17).Value = -0.5
$ws.Cells.Item(196, 18).Value = 1.85
$ws.Cells.Item(196, 19).Value = 1.95
$ws.Cells.Item(196, 20).Value = 2.25
$ws.Cells.Item(196, 21).Value = 1.975
$ws.Cells.Item(196, 22).Value = 1.825
$ws.Cells.Item(196, 23).Value = 0.909
$ws.Cells.Item(196, 25).Value = -1
$ws.Cells.Item(196, 26).Value = 0.8500000000000001
$ws.Cells.Item(196, 27).Value = -1
$ws.Cells.Item(196, 28).Value = 0.9750000000000001
$ws.Cells.Item(196, 29).Value = -1

# Row 197
$ws.Cells.Item(197, 2).Value = 6343950
$ws.Cells.Item(197, 6).Value = "Ituano"
$ws.Cells.Item(197, 7).Value = "Vitoria"
$ws.Cells.Item(197, 8).Value = 0
$ws.Cells.Item(197, 9).Value = 2
$ws.Cells.Item(197, 10).Value = "A"
$ws.Cells.Item(197, 11).Value = 2.5
$ws.Cells.Item(197, 12).Value = 2.9
$ws.Cells.Item(197, 13).Value = 2.9
$ws.Cells.Item(197, 14).Value = 2.3
$ws.Cells.Item(197, 15).Value = 3
$ws.Cells.Item(197, 16).Value = 3.4
$ws.Cells.Item(197, 17).Value = -0.25
$ws.Cells.Item(197, 18).Value = 1.9
$ws.Cells.Item(197, 19).Value = 1.9
$ws.Cells.Item(197, 20).Value = 2
$ws.Cells.Item(197, 21).Value = 1.875
$ws.Cells.Item(197, 22).Value = 1.925
$ws.Cells.Item(197, 23).Value = -1
$ws.Cells.Item(197, 25).Value = 2.4
$ws.Cells.Item(197, 26).Value = -1
$ws.Cells.Item(197, 27).Value = 0.8999999999999999
$ws.Cells.Item(197, 28).Value = 0
$ws.Cells.Item(197, 29).Value = -0

# Row 276
$ws.Cells.Item(276, 2).Value = 6673197
$ws.Cells.Item(276, 6).Value = "Vitoria"
$ws.Cells.Item(276, 7).Value = "Sport Recife"
$ws.Cells.Item(276, 8).Value = 1
$ws.Cells.Item(276, 9).Value = 0
$ws.Cells.Item(276, 12).Value = 3.2
$ws.Cells.Item(276, 13).Value = 4.333
$ws.Cells.Item(276, 14).Value = 2.9
$ws.Cells.Item(276, 15).Value = 3
$ws.Cells.Item(276, 16).Value = 2.625
$ws.Cells.Item(276, 17).Value = 0
$ws.Cells.Item(276, 18).Value = 2.05
$ws.Cells.Item(276, 19).Value = 1.8
$ws.Cells.Item(276, 21).Value = 1.9
$ws.Cells.Item(276, 22).Value = 1.95
$ws.Cells.Item(276, 23).Value = 1.9
$ws.Cells.Item(276, 26).Value = 1.05
$ws.Cells.Item(276, 28).Value = -1
$ws.Cells.Item(276, 29).Value = 0.95

# Row 277
$ws.Cells.Item(277, 2).Value = 6680024
$ws.Cells.Item(277, 6).Value = "Vila Nova"
$ws.Cells.Item(277, 7).Value = "Ceara"
$ws.Cells.Item(277, 8).Value = 3
$ws.Cells.Item(277, 9).Value = 1
$ws.Cells.Item(277, 10).Value = "H"
$ws.Cells.Item(277, 11).Value = 1.666
$ws.Cells.Item(277, 12).Value = 3.5
$ws.Cells.Item(277, 13).Value = 5.5
$ws.Cells.Item(277, 14).Value = 1.571
$ws.Cells.Item(277, 15).Value = 3.6
$ws.Cells.Item(277, 16).Value = 7.5
$ws.Cells.Item(277, 17).Value = -1
$ws.Cells.Item(277, 20).Value = 2.25
$ws.Cells.Item(277, 21).Value = 1.975
$ws.Cells.Item(277, 22).Value = 1.875
$ws.Cells.Item(277, 23).Value = 0.571
$ws.Cells.Item(277, 24).Value = -1
$ws.Cells.Item(277, 26).Value = 1.05
$ws.Cells.Item(277, 27).Value = -1
$ws.Cells.Item(277, 28).Value = 0.9750000000000001
$ws.Cells.Item(277, 29).Value = -1

# Row 279
$ws.Cells.Item(279, 2).Value = 6679229
$ws.Cells.Item(279, 6).Value = "EC Juventude"
$ws.Cells.Item(279, 7).Value = "Ponte Preta"
$ws.Cells.Item(279, 8).Value = 0
$ws.Cells.Item(279, 9).Value = 0
$ws.Cells.Item(279, 10).Value = "D"
$ws.Cells.Item(279, 11).Value = 1.45
$ws.Cells.Item(279, 12).Value = 4.5
$ws.Cells.Item(279, 13).Value = 6.5
$ws.Cells.Item(279, 14).Value = 1.8
$ws.Cells.Item(279, 15).Value = 3.3
$ws.Cells.Item(279, 16).Value = 5.25
$ws.Cells.Item(279, 17).Value = -0.75
$ws.Cells.Item(279, 20).Value = 2
$ws.Cells.Item(279, 21).Value = 2.1
$ws.Cells.Item(279, 22).Value = 1.775
$ws.Cells.Item(279, 23).Value = -1
$ws.Cells.Item(279, 24).Value = 2.3
$ws.Cells.Item(279, 26).Value = -1
$ws.Cells.Item(279, 27).Value = 0.8
$ws.Cells.Item(279, 28).Value = -1
$ws.Cells.Item(279, 29).Value = 0.7749999999999999

# Row 280
$ws.Cells.Item(280, 2).Value = 6673199
$ws.Cells.Item(280, 6).Value = "Mirassol"
$ws.Cells.Item(280, 7).Value = "Atletico GO"
$ws.Cells.Item(280, 8).Value = 4
$ws.Cells.Item(280, 9).Value = 1
$ws.Cells.Item(280, 12).Value = 3.4
$ws.Cells.Item(280, 13).Value = 4
$ws.Cells.Item(280, 14).Value = 2.2
$ws.Cells.Item(280, 15).Value = 3.25
$ws.Cells.Item(280, 16).Value = 3.5
$ws.Cells.Item(280, 17).Value = -0.25
$ws.Cells.Item(280, 18).Value = 1.925
$ws.Cells.Item(280, 19).Value = 1.925
$ws.Cells.Item(280, 21).Value = 1.85
$ws.Cells.Item(280, 22).Value = 2
$ws.Cells.Item(280, 23).Value = 1.2
$ws.Cells.Item(280, 26).Value = 0.925
$ws.Cells.Item(280, 28).Value = 0.8500000000000001
$ws.Cells.Item(280, 29).Value = -1

# Row 284
$ws.Cells.Item(284, 2).Value = 6693029
$ws.Cells.Item(284, 6).Value = "Atletico GO"
$ws.Cells.Item(284, 7).Value = "Guarani"
$ws.Cells.Item(284, 8).Value = 3
$ws.Cells.Item(284, 9).Value = 0
$ws.Cells.Item(284, 10).Value = "H"
$ws.Cells.Item(284, 11).Value = 1.45
$ws.Cells.Item(284, 12).Value = 4
$ws.Cells.Item(284, 13).Value = 8
$ws.Cells.Item(284, 14).Value = 1.333
$ws.Cells.Item(284, 15).Value = 4.8
$ws.Cells.Item(284, 16).Value = 10
$ws.Cells.Item(284, 17).Value = -1.25
$ws.Cells.Item(284, 18).Value = 1.8
$ws.Cells.Item(284, 19).Value = 2
$ws.Cells.Item(284, 20).Value = 2.5
$ws.Cells.Item(284, 21).Value = 1.95
$ws.Cells.Item(284, 22).Value = 1.85
$ws.Cells.Item(284, 23).Value = 0.333
$ws.Cells.Item(284, 25).Value = -1
$ws.Cells.Item(284, 26).Value = 0.8
$ws.Cells.Item(284, 27).Value = -1
$ws.Cells.Item(284, 28).Value = 0.95

# Row 285
$ws.Cells.Item(285, 2).Value = 6689429
$ws.Cells.Item(285, 6).Value = "ABC"
$ws.Cells.Item(285, 7).Value = "Vila Nova"
$ws.Cells.Item(285, 8).Value = 3
$ws.Cells.Item(285, 9).Value = 2
$ws.Cells.Item(285, 10).Value = "H"
$ws.Cells.Item(285, 11).Value = 8
$ws.Cells.Item(285, 12).Value = 4.75
$ws.Cells.Item(285, 13).Value = 1.363
$ws.Cells.Item(285, 14).Value = 6.5
$ws.Cells.Item(285, 15).Value = 4.2
$ws.Cells.Item(285, 16).Value = 1.45
$ws.Cells.Item(285, 17).Value = 1
$ws.Cells.Item(285, 18).Value = 2
$ws.Cells.Item(285, 19).Value = 1.8
$ws.Cells.Item(285, 20).Value = 2.25
$ws.Cells.Item(285, 21).Value = 1.9
$ws.Cells.Item(285, 22).Value = 1.9
$ws.Cells.Item(285, 23).Value = 5.5
$ws.Cells.Item(285, 25).Value = -1
$ws.Cells.Item(285, 26).Value = 1
$ws.Cells.Item(285, 27).Value = -1
$ws.Cells.Item(285, 28).Value = 0.8999999999999999
$ws.Cells.Item(285, 29).Value = -1

# Row 286
$ws.Cells.Item(286, 2).Value = 6693367
$ws.Cells.Item(286, 6).Value = "Sport Recife"
$ws.Cells.Item(286, 7).Value = "Sampaio Correa"
$ws.Cells.Item(286, 8).Value = 4
$ws.Cells.Item(286, 9).Value = 1
$ws.Cells.Item(286, 11).Value = 1.444
$ws.Cells.Item(286, 12).Value = 4
$ws.Cells.Item(286, 13).Value = 7
$ws.Cells.Item(286, 14).Value = 1.533
$ws.Cells.Item(286, 15).Value = 4
$ws.Cells.Item(286, 16).Value = 6
$ws.Cells.Item(286, 17).Value = -1
$ws.Cells.Item(286, 18).Value = 1.975
$ws.Cells.Item(286, 19).Value = 1.825
$ws.Cells.Item(286, 21).Value = 1.775
$ws.Cells.Item(286, 22).Value = 2.025
$ws.Cells.Item(286, 23).Value = 0.5329999999999999
$ws.Cells.Item(286, 26).Value = 0.9750000000000001
$ws.Cells.Item(286, 28).Value = 0.7749999999999999

# Row 287
$ws.Cells.Item(287, 2).Value = 6689350
$ws.Cells.Item(287, 6).Value = "Tombense MG"
$ws.Cells.Item(287, 7).Value = "Mirassol"
$ws.Cells.Item(287, 8).Value = 0
$ws.Cells.Item(287, 9).Value = 1
$ws.Cells.Item(287, 10).Value = "A"
$ws.Cells.Item(287, 11).Value = 3.2
$ws.Cells.Item(287, 12).Value = 3
$ws.Cells.Item(287, 13).Value = 2.25
$ws.Cells.Item(287, 14).Value = 3
$ws.Cells.Item(287, 15).Value = 3.25
$ws.Cells.Item(287, 16).Value = 2.3
$ws.Cells.Item(287, 17).Value = 0.25
$ws.Cells.Item(287, 18).Value = 1.775
$ws.Cells.Item(287, 19).Value = 2.025
$ws.Cells.Item(287, 21).Value = 2
$ws.Cells.Item(287, 22).Value = 1.8
$ws.Cells.Item(287, 23).Value = -1
$ws.Cells.Item(287, 25).Value = 1.3
$ws.Cells.Item(287, 26).Value = -1
$ws.Cells.Item(287, 27).Value = 1.025
$ws.Cells.Item(287, 28).Value = -1
$ws.Cells.Item(287, 29).Value = 0.8

# Row 288
$ws.Cells.Item(288, 2).Value = 6693030
$ws.Cells.Item(288, 6).Value = "Ceara"
$ws.Cells.Item(288, 7).Value = "EC Juventude"
$ws.Cells.Item(288, 8).Value = 1
$ws.Cells.Item(288, 9).Value = 3
$ws.Cells.Item(288, 10).Value = "A"
$ws.Cells.Item(288, 11).Value = 3.25
$ws.Cells.Item(288, 12).Value = 3.4
$ws.Cells.Item(288, 13).Value = 2.1
$ws.Cells.Item(288, 14).Value = 3.3
$ws.Cells.Item(288, 15).Value = 3.2
$ws.Cells.Item(288, 16).Value = 2.3
$ws.Cells.Item(288, 17).Value = 0.25
$ws.Cells.Item(288, 18).Value = 1.875
$ws.Cells.Item(288, 19).Value = 1.975
$ws.Cells.Item(288, 21).Value = 1.825
$ws.Cells.Item(288, 22).Value = 2.025
$ws.Cells.Item(288, 23).Value = -1
$ws.Cells.Item(288, 25).Value = 1.3
$ws.Cells.Item(288, 26).Value = -1
$ws.Cells.Item(288, 27).Value = 0.9750000000000001
$ws.Cells.Item(288, 28).Value = 0.825

# Row 289
$ws.Cells.Item(289, 2).Value = 6693031
$ws.Cells.Item(289, 6).Value = "Chapecoense"
$ws.Cells.Item(289, 7).Value = "Vitoria"
$ws.Cells.Item(289, 8).Value = 3
$ws.Cells.Item(289, 11).Value = 1.8
$ws.Cells.Item(289, 12).Value = 3.6
$ws.Cells.Item(289, 13).Value = 4.333
$ws.Cells.Item(289, 14).Value = 1.615
$ws.Cells.Item(289, 16).Value = 5.25
$ws.Cells.Item(289, 18).Value = 2.025
$ws.Cells.Item(289, 19).Value = 1.775
$ws.Cells.Item(289, 23).Value = 0.615
$ws.Cells.Item(289, 26).Value = 1.025

# Row 290
$ws.Cells.Item(290, 2).Value = 6693028
$ws.Cells.Item(290, 6).Value = "Ponte Preta"
$ws.Cells.Item(290, 7).Value = "CRB"
$ws.Cells.Item(290, 9).Value = 0
$ws.Cells.Item(290, 11).Value = 1.727
$ws.Cells.Item(290, 12).Value = 3.5
$ws.Cells.Item(290, 13).Value = 4
$ws.Cells.Item(290, 14).Value = 1.7
$ws.Cells.Item(290, 15).Value = 3.6
$ws.Cells.Item(290, 16).Value = 5
$ws.Cells.Item(290, 17).Value = -0.75
$ws.Cells.Item(290, 18).Value = 1.975
$ws.Cells.Item(290, 19).Value = 1.875
$ws.Cells.Item(290, 20).Value = 2
$ws.Cells.Item(290, 22).Value = 2.1
$ws.Cells.Item(290, 23).Value = 0.7
$ws.Cells.Item(290, 26).Value = 0.9750000000000001

# Row 292
$ws.Cells.Item(292, 2).Value = 6689428
$ws.Cells.Item(292, 6).Value = "Botafogo SP"
$ws.Cells.Item(292, 7).Value = "Londrina"
$ws.Cells.Item(292, 9).Value = 1
$ws.Cells.Item(292, 10).Value = "A"
$ws.Cells.Item(292, 11).Value = 1.615
$ws.Cells.Item(292, 12).Value = 3.5
$ws.Cells.Item(292, 13).Value = 6
$ws.Cells.Item(292, 14).Value = 1.615
$ws.Cells.Item(292, 15).Value = 3.8
$ws.Cells.Item(292, 16).Value = 5.75
$ws.Cells.Item(292, 17).Value = -0.75
$ws.Cells.Item(292, 18).Value = 1.825
$ws.Cells.Item(292, 19).Value = 1.975
$ws.Cells.Item(292, 20).Value = 2.5
$ws.Cells.Item(292, 21).Value = 2
$ws.Cells.Item(292, 22).Value = 1.8
$ws.Cells.Item(292, 24).Value = -1
$ws.Cells.Item(292, 25).Value = 4.75
$ws.Cells.Item(292, 26).Value = -1
$ws.Cells.Item(292, 27).Value = 0.9750000000000001
$ws.Cells.Item(292, 29).Value = 0.8

# Row 293
$ws.Cells.Item(293, 2).Value = 6689427
$ws.Cells.Item(293, 6).Value = "Avai"
$ws.Cells.Item(293, 7).Value = "Ituano"
$ws.Cells.Item(293, 9).Value = 0
$ws.Cells.Item(293, 10).Value = "D"
$ws.Cells.Item(293, 11).Value = 1.95
$ws.Cells.Item(293, 12).Value = 3.1
$ws.Cells.Item(293, 13).Value = 4.2
$ws.Cells.Item(293, 14).Value = 2.4
$ws.Cells.Item(293, 15).Value = 3.2
$ws.Cells.Item(293, 16).Value = 3.2
$ws.Cells.Item(293, 17).Value = -0.25
$ws.Cells.Item(293, 18).Value = 2.05
$ws.Cells.Item(293, 19).Value = 1.75
$ws.Cells.Item(293, 20).Value = 2.25
$ws.Cells.Item(293, 21).Value = 1.825
$ws.Cells.Item(293, 22).Value = 1.975
$ws.Cells.Item(293, 24).Value = 2.2
$ws.Cells.Item(293, 25).Value = -1
$ws.Cells.Item(293, 26).Value = -0.5
$ws.Cells.Item(293, 27).Value = 0.375
$ws.Cells.Item(293, 29).Value = 0.9750000000000001
